# Added exception handling for when plantcode is not found in locations.json
# The example/placeholder Shipper Plantcode on the Input sheet is updated
# from the old sentinel "A999" to "A001".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")

# Update the Shipper Plantcode sample value (B6) from "A999" to "A001"
$ws.Range("B6").Value = "A001"

# Move the active selection to B8, matching the saved cursor position
$ws.Activate()
$ws.Range("B8").Select()
